$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.216.16"
$ws.Range("E2").Value = "  -3.88%  "

$ws.Range("D3").Value = "1.659.61"
$ws.Range("E3").Value = "  -2.51%  "

$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.12"
$ws.Range("E5").Value = "  -2.37%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5156"
$ws.Range("E6").Value = "  -2.73%  "

$ws.Range("E7").Value = "  +0.25%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2582"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06451"
$ws.Range("E9").Value = "  -1.98%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.01"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07818"
$ws.Range("E11").Value = "  +2.63%  "

$ws.Range("D12").Value = "1.668.26"
$ws.Range("E12").Value = "  -1.93%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.305"
$ws.Range("E13").Value = "  -4.25%  "

$ws.Range("D14").Value = "1.887.33"
$ws.Range("E14").Value = "  -2.54%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5553"
$ws.Range("E15").Value = "  -4.01%  "

$ws.Range("D16").Value = "0.0₅8076"
$ws.Range("E16").Value = "  -0.79%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.38"
$ws.Range("E17").Value = "  -4.46%  "

$ws.Range("D18").Value = "26.239.13"
$ws.Range("E18").Value = "  -3.78%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "211.72"
$ws.Range("E19").Value = "  -1.73%  "

$ws.Range("E20").Value = "  +0.27%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.08"
$ws.Range("E22").Value = "  -2.40%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.995"
$ws.Range("E23").Value = "  +0.48%  "

$ws.Range("E24").Value = "  +0.17%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.58"
$ws.Range("E25").Value = "  +0.30%  "

$ws.Range("E26").Value = "  +3.07%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1170"
$ws.Range("E27").Value = "  -2.19%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.991"
$ws.Range("E28").Value = "  -2.88%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.81"
$ws.Range("E29").Value = "  -2.09%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05221"
$ws.Range("E30").Value = "  -2.63%  "

$ws.Range("E31").Value = "  -2.24%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.368"
$ws.Range("E32").Value = "  -2.79%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.230"
$ws.Range("E33").Value = "  -4.95%  "

$ws.Range("E34").Value = "  -3.96%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.760"
$ws.Range("E35").Value = "  -3.50%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9338"
$ws.Range("E36").Value = "  -1.22%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.370"
$ws.Range("E37").Value = "  -1.87%  "

$ws.Range("D38").Value = "1.177.06"
$ws.Range("E38").Value = "  +13.22%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5706"
$ws.Range("E39").Value = "  -1.89%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01595"
$ws.Range("E40").Value = "  -2.04%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8453"
$ws.Range("E41").Value = "  +0.82%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.004"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.683"
$ws.Range("E43").Value = "  -1.56%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.70"
$ws.Range("E44").Value = "  -0.14%  "

$ws.Range("D45").Value = "1.797.52"
$ws.Range("E45").Value = "  -2.55%  "

$ws.Range("E46").Value = "  -1.33%  "

$ws.Range("E47").Value = "  +0.46%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "56.07"
$ws.Range("E48").Value = "  -2.93%  "

$ws.Range("E49").Value = "  -0.20%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.882"
$ws.Range("E50").Value = "  -2.30%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05063"
$ws.Range("E51").Value = "  -3.08%  "

